# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts.
#
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a
# set of rows on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 14;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 19;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 22;  Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 27;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 52;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 66;  Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 74;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 78;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 85;  Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 100; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 110; Tag = "%";  Label = "Uninterpretable" },
    @{ Row = 118; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 124; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 143; Tag = "ba"; Label = "Appreciation" },
    @{ Row = 144; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 148; Tag = "aa"; Label = "Agree/Accept" },
    @{ Row = 157; Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 160; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 165; Tag = "%";  Label = "Uninterpretable" },
    @{ Row = 167; Tag = "aa"; Label = "Agree/Accept" },
    @{ Row = 176; Tag = "%";  Label = "Uninterpretable" },
    @{ Row = 201; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 211; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 214; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 225; Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 226; Tag = "ba"; Label = "Appreciation" },
    @{ Row = 254; Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 255; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 257; Tag = "sd"; Label = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Label
}
